# Heal QA Automation test_data.xlsx update
# - Added state and phone to Account sheet test data
# - Added member id / group id / payer / phone to Patient sheet test data
# - Updated ZipCode sample values on both sheets
# - Widened the Value column on the Patient sheet
# - Updated active-cell selections left by the editing session

$wb = $excel.ActiveWorkbook

$wsAccount = $wb.Worksheets.Item("Account")
$wsPatient = $wb.Worksheets.Item("Patient")

# --- Patient sheet ("Patient") -------------------------------------------------
# Existing sample zip code value changed
$wsPatient.Range("B7").Value = 90210

# New rows: MemberId / GroupId keys first, then their values, then Payer pair,
# then the Phone row (matches the order the values were typed in originally).
$wsPatient.Range("A11").Value = "MemberId"
$wsPatient.Range("A12").Value = "GroupId"
$wsPatient.Range("B11").Value = "COST_ESTIMATES_025"
$wsPatient.Range("B12").Value = "BC001"
$wsPatient.Range("A13").Value = "Payer"
$wsPatient.Range("B13").Value = "Anthem Blue Cross"
$wsPatient.Range("A14").Value = "Phone"
$wsPatient.Range("B14").Value = 2015555555

# Widen column B (Value) so the new longer values are readable
$wsPatient.Columns.Item(2).ColumnWidth = 20

# --- Account sheet ("Account") --------------------------------------------------
# Existing sample zip code value changed
$wsAccount.Range("B7").Value = 90210

# New rows: Phone and State
$wsAccount.Range("A15").Value = "Phone"
$wsAccount.Range("B15").Value = 2015555555
$wsAccount.Range("A16").Value = "State"
$wsAccount.Range("B16").Value = "California"

# --- Leave the workbook with the selections/active sheet from the edit session --
[void]$wsPatient.Activate()
[void]$wsPatient.Range("B17").Select()

[void]$wsAccount.Activate()
[void]$wsAccount.Range("D13").Select()
